# Applies the attendance_reports sync update to the Session Analysis Results sheet.

# Helper: write a literal text value into a cell without Excel's
# "looks like a percentage" auto-conversion to a number (which would
# also silently swap in a new percentage number format / style).
# We do this by writing a text-returning formula (T() forces text) and
# then immediately collapsing it to a static value via PasteSpecial
# values-only, which leaves the original style/number format untouched.
function Set-LiteralText {
    param($ws, $addr, [string]$text)
    $cell = $ws.Range($addr)
    $escaped = $text -replace '"', '""'
    $cell.Formula = '=T("' + $escaped + '")'
    $cell.Copy()
    $cell.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reordered "Recorded By" email lists (content unchanged, order shuffled) ---

$ws.Range("G2").Value  = "nourhan.mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg"
$ws.Range("G24").Value = "nourhan.mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg"

$ws.Range("G18").Value = "shorokmohamed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg"
$ws.Range("G40").Value = "shorokmohamed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg"

$ws.Range("G19").Value  = "wafaa.ebida@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg"
$ws.Range("G41").Value  = "wafaa.ebida@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg"
$ws.Range("G150").Value = "wafaa.ebida@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg"
$ws.Range("G172").Value = "wafaa.ebida@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg"

$ws.Range("G46").Value = "hend_mahmoud@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg"
$ws.Range("G68").Value = "hend_mahmoud@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg"

$ws.Range("G54").Value  = "maimustafa@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, merna.said@med.asu.edu.eg"
$ws.Range("G76").Value  = "maimustafa@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, merna.said@med.asu.edu.eg"
$ws.Range("G98").Value  = "maimustafa@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, merna.said@med.asu.edu.eg"
$ws.Range("G120").Value = "maimustafa@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, merna.said@med.asu.edu.eg"

$ws.Range("G62").Value = "wafaa.ebida@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg"
$ws.Range("G84").Value = "wafaa.ebida@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg"

$ws.Range("G96").Value  = "aml.awwad@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg"
$ws.Range("G118").Value = "aml.awwad@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg"

$ws.Range("G106").Value = "wafaa.ebida@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg"
$ws.Range("G128").Value = "wafaa.ebida@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg"

$ws.Range("G107").Value = "wafaa.ebida@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg"
$ws.Range("G129").Value = "wafaa.ebida@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg"

$ws.Range("G112").Value = "Mohammedeltanany@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg"

$ws.Range("G134").Value = "hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"

$ws.Range("G142").Value = "esraa.mostafa@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, merna.said@med.asu.edu.eg"
$ws.Range("G164").Value = "esraa.mostafa@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, merna.said@med.asu.edu.eg"

$ws.Range("G156").Value = "majorelle.magdy@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, Mohammedeltanany@med.asu.edu.eg"

# --- Class Statistics summary box (K2 block) ---
$ws.Range("L6").Value = 41
$ws.Range("L7").Value = 11
Set-LiteralText $ws "L9"  "23.3%"
Set-LiteralText $ws "L10" "32.8%"

# --- Per-group statistics rows (K16:S16 = group A2, K19:S19 = group B1) ---
$ws.Range("O16").Value = 5
$ws.Range("P16").Value = 2
Set-LiteralText $ws "R16" "22.7%"
Set-LiteralText $ws "S16" "31.6%"

$ws.Range("O19").Value = 6
$ws.Range("P19").Value = 2
Set-LiteralText $ws "R19" "27.3%"
Set-LiteralText $ws "S19" "31.5%"

# --- Rows that flipped from "Not Recorded" to "Recorded" (style + content) ---
# Row 32: Year 2 / A2 / MICROBIOLOGY / session 1
$src = $ws.Range("A24:I24")
$src.Copy()
$ws.Range("A32:I32").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("G32").Value = "amany.raafat@med.asu.edu.eg"
$ws.Range("H32").Value = "20/217"
$ws.Range("I32").Value = "Recorded"

# Row 90: Year 2 / B1 / ANATOMY / session 1
$src = $ws.Range("A2:I2")
$src.Copy()
$ws.Range("A90:I90").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("G90").Value = "Mohammedeltanany@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg"
$ws.Range("H90").Value = "34/154"
$ws.Range("I90").Value = "Recorded"

# Row 91: Year 2 / B1 / ANATOMY / session 2
$src = $ws.Range("A2:I2")
$src.Copy()
$ws.Range("A91:I91").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("G91").Value = "Mohammedeltanany@med.asu.edu.eg"
$ws.Range("H91").Value = "50/154"
$ws.Range("I91").Value = "Recorded"

$excel.CutCopyMode = $false
